# Price-tracker update for 2026-02-07.
# Appends one new row (row 38) to the single data sheet with the freshly
# scraped Date/Price, mirroring the existing rows (Discount=0, Incredible=0
# i.e. shared-string "0", same as every other "in stock / no discount" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 37
$newRow  = 38

# 1) Clone the previous row verbatim (values + "plain text" cell typing) so
#    the new row starts out identical to row 37 - this avoids Excel's
#    automatic number/date reinterpretation that a fresh .Value= assignment
#    would trigger for numeric-looking / date-looking text.
$srcRange = $ws.Range("A$lastRow" + ":D$lastRow")
$dstRange = $ws.Range("A$newRow" + ":D$newRow")
$srcRange.Copy($dstRange)

# 2) Overwrite Date (col A) and Price (col B) with the new scraped values.
#    Routing the literal text through a formula ("=""text""") and then
#    collapsing it back down to a plain value via Copy + PasteSpecial
#    (xlPasteValues) writes it as ordinary text without Excel "helpfully"
#    turning "2026-02-07" into a date serial or "759000" into a number.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Formula = "=""2026-02-07"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$priceCell = $ws.Cells.Item($newRow, 2)
$priceCell.Formula = "=""759000"""
$priceCell.Copy()
$priceCell.PasteSpecial(-4163)

# Columns C (Discount) and D (Incredible) stay "0"/"0", already copied from
# row 37 above, so nothing else to do.

$excel.CutCopyMode = 0
